$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hour values for row 2 (week commencing 43122)
$ws.Range("B2").Value = 3.75
$ws.Range("C2").Value = 2.25
$ws.Range("E2").Value = 5.25
$ws.Range("F2").Value = 5.5
$ws.Range("G2").Value = 6.5
$ws.Range("H2").Value = 7.5

# Update hour values for row 3 (week commencing 43129)
$ws.Range("B3").Value = 6.5
$ws.Range("E3").Value = 5.25
$ws.Range("F3").Value = 5.5
$ws.Range("G3").Value = 8.5
$ws.Range("H3").Value = 7.5

# Update hour values for row 4 (week commencing 43136)
$ws.Range("B4").Value = 7.5
$ws.Range("C4").Value = 1.25
$ws.Range("E4").Value = 5.75
$ws.Range("F4").Value = 7.25

# Update the active cell selection
$ws.Range("N8").Select()
